$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("I2").Value = -1.724731182795699
$ws.Range("J2").Value = -1.682100115424336
$ws.Range("K2").Value = 0.397
$ws.Range("L2").Value = 0.4268817204301075
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("U2").Value = 5.5
$ws.Range("V2").Value = 0.05191617896922786
$ws.Range("X2").Value = 0.03846338425444983
$ws.Range("Z2").Value = 0.1604831751509923
$ws.Range("AB2").Value = 0.03846338425444983
$ws.Range("AD2").Value = 0.012
$ws.Range("AF2").Value = 0.012
$ws.Range("AG2").Value = -5.488
$ws.Range("AH2").Value = 0.0001132588341890667
$ws.Range("AI2").Value = 0.0006684491978609625
$ws.Range("AJ2").Value = -0.05463305857523992
$ws.Range("AK2").Value = -0.4407324124638612
$ws.Range("AM2").Value = -2.465
$ws.Range("W2").Value = 0.06723076923076923
$ws.Range("Y2").Value = 0.02876738497631941
$ws.Range("AA2").Value = 0.08027579112373862
$ws.Range("AC2").Value = 0.04174669577572879
$ws.Range("AQ2").Value = 0.6507099391480731
$ws.Range("T2").ClearContents()

# --- Row 3 updates ---
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "AIQ Limited (LSE:AIQ)"
$ws.Range("I3").Value = -21.41935483870968
$ws.Range("J3").Value = -21.41935483870968
$ws.Range("K3").Value = -1.09
$ws.Range("L3").Value = -35.16129032258065
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 4.05
$ws.Range("V3").Value = 0.27
$ws.Range("X3").Value = 0.03846338425444983
$ws.Range("AB3").Value = 0.03846338425444983
$ws.Range("AD3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -4.05
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.3698630136986302
$ws.Range("AK3").Value = -1.919431279620853
$ws.Range("AM3").Value = -0.025
$ws.Range("W3").Value = -0.2416851441241686
$ws.Range("Y3").Value = -0.2801485283786184
$ws.Range("Z3").Value = -0.0574074074074074
$ws.Range("AA3").Value = 1.22962962962963
$ws.Range("AC3").Value = 1.19116624537518
$ws.Range("AQ3").Value = 26.56

# --- Row 4 updates ---
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "Mineral & Financial Investments Limited (AIM:MAFL)"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.6006674082313682
$ws.Range("J4").Value = 0.5561263927904828
$ws.Range("K4").Value = 0.437
$ws.Range("L4").Value = 0.4860956618464961
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0.34
$ws.Range("V4").Value = 0.07172995780590717
$ws.Range("W4").Value = 0.06723076923076923
$ws.Range("X4").Value = 0.03853721560699273
$ws.Range("Y4").Value = 0.0286935536237765
$ws.Range("Z4").Value = 0.1443481053307643
$ws.Range("AA4").Value = 0.08027579112373862
$ws.Range("AB4").Value = 0.03852909534800982
$ws.Range("AC4").Value = 0.04174669577572879
$ws.Range("AD4").Value = 0.012
$ws.Range("AF4").Value = 0.012
$ws.Range("AG4").Value = -0.328
$ws.Range("AH4").Value = 0.002525252525252525
$ws.Range("AI4").Value = 0.00176678445229682
$ws.Range("AJ4").Value = -0.0743427017225748
$ws.Range("AK4").Value = -0.05083694978301302
$ws.Range("T4").ClearContents()

# --- Row 5 (new row) ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "Cayman Islands"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "Alussa Energy Acquisition Corp. (NYSE:ALUS)"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "Investments & Asset Management"
$ws.Range("K5").Value = 1.05
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 1.11
$ws.Range("V5").Value = 0.01287703016241299
$ws.Range("W5").Value = 55.26315789473684
$ws.Range("X5").Value = 0.03846338425444983
$ws.Range("Y5").Value = 55.22469451048239
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = -13.83177570093458
$ws.Range("AB5").Value = 0.03846338425444983
$ws.Range("AC5").Value = -13.87023908518903
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = -1.11
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = -0.0130450111646492
$ws.Range("AK5").Value = -0.2853470437017995
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = -2.44
$ws.Range("AQ5").Value = 0.6065573770491803
